$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NAME")
Write-Host "Sheet name: $($ws.Name)"
